# Clear the stray placeholder zero values on the "Solvent Properties" sheet.
# These cells originally held a literal numeric 0 where no real measurement
# exists; the updated data-export script now leaves them blank (an empty
# text value) instead of writing a numeric 0.
#
# A plain empty-string assignment collapses the cell entirely (no cell left
# behind at all), so instead we type a lone apostrophe - Excel's "force
# text" prefix - which commits an actual empty *text* cell, then strip the
# leftover quote-prefix formatting so no visual/style residue remains.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Solvent Properties")

$cellsToClear = @(
    "H5", "J5",
    "H9", "I9", "J9",
    "H11", "I11", "J11",
    "J13",
    "J14",
    "J15",
    "J17",
    "I18", "J18",
    "J22",
    "J23",
    "J24",
    "I26", "J26",
    "I27", "J27",
    "I28", "J28",
    "I29",
    "J30",
    "J31",
    "J32",
    "J33",
    "J34",
    "H36", "I36", "J36",
    "J38",
    "J39",
    "J40",
    "E42", "H42", "I42", "J42",
    "H44", "I44", "J44",
    "J45",
    "J47",
    "J50",
    "J51",
    "J52"
)

foreach ($addr in $cellsToClear) {
    $cell = $ws.Range($addr)
    $cell.Value = "'"
    $cell.ClearFormats()
}
